$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$i21Text = @'
dynex capital dx a monthly dividend stock under with big yield dividend stock in the retail sector under guess ge stand out in trump medium stock surge on bitcoin stockpile los angeles time owner reveals he taking newspaper public flower food flo a reliable dividend stock still trading under top meme coin list flipped balloon cat presale blast past k at stage universal music group file for u public offering los angeles time to become publicly traded company variety energy transfer et a rare dividend stock trading below keycorps key dividend history what make it stand out under how europe trade bazooka could be a last resort against trump tariff best ai accounting software quickbooks recognized a top ai accounting solution in report by better business advice break the link between speed and success on your team general motor is set to report earnings before the bell here what wall street expects cnbc tripadvisor reveals this year toprated travel experience worldwide with traveler choice award best of the best thing to do gabbards russian interference claim directly contradict what other trump official have said jpmorgan explores lending against client cryptocurrency financial time china problem with competition there too much of it the new york time top underrated apps to download for k smart tv hay testing more important this year baptist health announces new name for medical group ai adoption is growing but why secure scaling is lagging is exxons xom dividend safe and worth it for halal investor ev share in china june sale report hyundai motor robotics lab triumph at red dot award design concept with visionary innovation clout war jensen huang eclipse elon musk and tim cook in washington kia design triumph at red dot award design concept with award top honor for pv wkndr cocacolas ko global brand and dividend reliability for ethical investor france increase in nuclear hydropower in led to more electricity export is chevron cvx a good longterm stock for halal portfolio report goldman sachs advising bnsf on potential merger what halal investor should know about pepsico pep china begin construction of world largest hydropower dam is texas instrument txn a halal dividend stock for techsavvy investor philippine buyer unsure whether they can absorb tariff cost say trade association it wa not justice breonna taylor mother speaks out on officer sentencing opendoor verizon gamestop alphabet and trump medium why these stock are on investor radar today medtronics mdt role in medical device a unique fit for halal portfolio why target tgt appeal to shariahconscious dividend investor is kenvue kvue a halal dividend stock worth adding to your portfolio treasury yield creep higher a investor await fed chair powell speech over million pool sold across the u canada under recall after report of nine death cnn evaluating syscos syy financial health for halal dividend investing why kimberlyclarks kmb dividend stability appeal to halal investor taiwan hold unprecedented military drill cocacola is about to report earnings here what to expect over million pool sold across the u canada under recall after report of nine death shareholder investigation halper sadeh llc investigates sptn turn strm on behalf of shareholder why archerdanielsmidland is a strong halal dividend stock choice gen z is returning to religion why shareholder investigation halper sadeh llc investigates vbtx hmst vrna str on behalf of shareholder regal rexnord corporation declares quarterly dividend of per share regal rexnord corporation declares quarterly dividend of per share ap business summarybrief at pm edt ap business summarybrief at pm edt ap business summarybrief at pm edt ap business summarybrief at pm edt ap technology summarybrief at pm edt ap technology summarybrief at pm edt kevin oleary epstein focus is poop on a stick looming over two case threatening musk car company is a single question can he be trusted ap news mindblowing milky way image captured by nasa telescope astrazeneca to invest billion in the u a pharma tariff weigh aca marketplace health plan enrollee could face subsidy cliff in here how to avoid it the real reason why innout is leaving california deseret news report that american man his family slaughtered by usbacked syrian regime gta member claim your free bravado buffalo evx today natural wine bar coming to schenectady pintsized ja solar and masdar deepen partnership with landmark gw pv project in azerbaijan mckinney to vote on funding for m jw marriott resort at craig ranch alaska airline resume flight after equipment failure at a data center ground all it plane new saratoga spring cannabis dispensary proposed for west side asia stock muted japan whipsaw after pm ishiba loses upper house investingcom judge partially block trump administration from enforcing funding ban against planned parenthood faraday future partner with habittrade to accelerate it entry into the web financial ecosystem denver mayor mike johnston call for learned hopefulness on homelessness housing other challenge bethel appoints dean for new unified business school peter thiel effect cathie wood pours million into bitmine immersion stock a founder fund backing ethereum spike fuel buzz fda approved moderna covid shot for kid behind kennedy back two hhs aide fired cnbc daily open investor look away from the trump administration call to relook the fed kate spade outlet shopper love this sleek and stylish tote bag thats on sale for only tesla retrofuturistic diner officially open a elon musk hint at more location techcrunch china affluent are a pessimistic about the economy a they were during the covid pandemic bird and segway partner to launch nextgen shared micromobility fleet across north america extrump advisor judy shelton slam federal reserve billion loss say violates the norm to have a fed chair who cant be fired hanwha philly shipyard land first u lng carrier order in decade from hanwha affiliate the decline of u housing affordability oil price dip a useu trade war loom ap business summarybrief at pm edt ap business summarybrief at pm edt ap business summarybrief at pm edt chinabacked hacker used microsoft flaw in attack defender say the washington post how europe trade bazooka could be a last resort against trump tariff astrazeneca to invest billion in the u a pharma tariff weigh kioxia announces industry first tb nvme ssd built for the demand of generative ai environment clout war jensen huang eclipse elon musk and tim cook in washington ted lasso return from across the pond for season with production now underway lutnick say trump will renegotiate trade pact with canada mexico next year farmer market moving out of downtown gilbert forget egg beef price are soaring u and uk aircraft carrier team up in the pacific cal thomas latenight tv not like it used to be the f jpo ha a new chief the u air force is saying goodbye to the a warthog trump and philippine leader plan to talk tariff and china at the white house russia is facing challenge with it fifthgeneration fighter letter new list same ol story dear neighbor living under a false mandate gabbards russian interference claim directly contradict what other trump official have said bryce fiedler it time to power down the chromebook experiment letter just say no regina dancy what will you choose chaos or community letter biking is a winwin for health and livability letter there are better way to spend billion than on a sport arena letter pension help nh retain the best police officer letter let our library be neutral the uk f stealth jet is finally flying back to britain there are many illegal marijuana farm but federal agent targeted california biggest legal one there are many illegal marijuana farm but federal agent targeted california biggest legal one npr got what it deserved public radio ceo denounced for defiant remark about nonwhite voice after congress killed federal funding septum urge philadelphia student to plan ahead due to possible service cut it wa not justice breonna taylor mother speaks out on officer sentencing iowa new summer food program a heavy lift for food pantry grand fork city council move forward with county district negotiation on apartment tax incentive go fk yourself colbert responds to trump post steering committee tap garbarino to serve a homeland security chair taiwan hold unprecedented military drill activist judge let mexican who allegedly decapitated a woman walk free jonah goldberg yes mr president theyre still talking about jeffrey epstein ted talley pro mundo an overreach a definition of baylor education over million pool sold across the u canada under recall after report of nine death sen john hoffman talk recovery in first interview since lawmaker attack dump truck roll over during crash on i in chelmsford can the coalition start parliament on a positive note fiona katauskas kevin oleary epstein focus is poop on a stick house speaker say there will be no vote on releasing epstein file before august recess preliminary finding to be announced tuesday from probe into fall river assisted living fire letter dont be late to a bok center show los angeles should follow new york lead and toll entry into the downtown area los angeles should follow new york lead and toll entry into the downtown area rep andrew garbarino to chair homeland security panel spearheading gop immigration agenda hunter biden claim father joe wa on ambien during disastrous debate against trump judge partially block trump administration from enforcing funding ban against planned parenthood schenectady police to place drone central park middle school roof day one of special session reveals rift over priority phillies get walkoff win over boston on catcher interference call with base loaded in th beach reopens after aggressive shark sighting near the redondo pier rep andrew garbarino win key vote in race for homeland security committee chair in cringe video zohran mamdani say he returning to africa but unfortunately coming back labor department proposes more than rule change in a push to deregulate workplace california ag bonta head coalition opposing usda plan to share snap participant data texas democrat vow not to make it easy for republican redistricting effort before midterm burdoin fire update july how to revitalize the abraham accord angel give up late run lead in loss to mets william harrys cousin found dead with gun nearby medicaid take center stage at fedorchaks first virtual town hall since big beautiful bill passage ted lasso return from across the pond for season with production now underway ice raid are hurting critical economic sector north orange county official gather to say one dead in jet crash near granite fall minnesota prosecutor want former oc supervisor to pay in restitution far below what county is seeking california pass major ceqa reform key takeaway for the regulated community jim vandehei rise of super journalist county commissioner eye nevada governor run in government asks disgraced former orange county supervisor to pay back nearly aocs campaign office vandalized with antiisrael message attorney general kwame raoul fight to preserve tps for immigrant joining coalition in legal battle against dhs plan house grind to a halt a gop try to shut down epstein vote hennepin county sheriff office boost community engagement through charitable educational and cultural initiative hhs say it will begin reform of organ donation system after federal investigation find horrifying problem federal judge block trump administration from defunding planned parenthood letter state park should have restaurant letter broken system leaf some tough choice letter pb is national treasure that deserves funding trump nominates former fox news contributor sara carter a drug czar cbss longawaited colbert cancellation letter to the editor july president trump letter to grieving mother of slain intern showcase his deep compassion and humanity tulsi gabbards file release raise question about claim from obama deputy national security advisor hundred of criminal case could soon be dismissed amid public defender work stoppage pentagon withdraws marine from los angeles day of texas special session start with democrat signaling they could flee over redistricting broken promise fractured friendship force a multimillion dollar lawsuit involving mlb superstar vladimir guerrero jr and singer marc anthony company florida law firm lesser landy smith siegel representing plaintiff jesse guerrero lawsuit challenge restriction on head start for kid in the u illegally fda name former pharmaceutical company executive to oversee u drug program martin luther king jr fbi record released by trump administration al sharpton say desperate attempt to distract people from epstein firestorm shareholder investigation halper sadeh llc investigates itos and zimv on behalf of shareholder cpr training eyed a high school grad requirement lawmaker hear gender notification bill trump labor department seek to deregulate workplace dismantling key protection for worker the golden dome ha a new czar epstein furor upends house for second week in row epstein furor upends house for second week in row epstein furor upends house for second week in row epstein furor upends house for second week in row epstein furor upends house for second week in row epstein furor upends house for second week in row tuesday brings more comfortable air but a big warm up is ahead justin baldonis insurer sue him wayfarer saying it policy dont cover blake lively harassment suit republican minnesota sen bruce anderson dy unexpectedly at age trump administration release record on martin luther king jr assassination danaher announces transition plan for chief financial officer we check on our favorite bond fund amid tariff volatility earnings preview what to expect from consolidated edison report westlake announces second quarter earnings conference call oscar health announces preliminary financial result for second quarter and revise guidance target is ending a beloved perk that shopper have relied on for year njcom global defense contractor burke product selects datavault ai for enhanced national defense and quavo fraud dispute secures million growth investment from spectrum equity highyield saving account rate today july rate are steady humana accelerates effort to eliminate prior authorization requirement to ensure a faster more seamless process standard chartered and village capital celebrate startup of inaugural woman in tech accelerator trump ai vision take shape a oracle and openai expand massive stargate infrastructure project clear to announce second quarter financial result on august heloc rate today july the home equity line of credit rate hang steady canada tariff could add to the cost of building a home by report warns it just wildly misleading why the administration latest allegation about the russia investigation dont add up the trump administration claim no one ha died due to u aid cut our trip to afghanistan suggests otherwise palm beach real estate mystery unfolds a buyer acquires property worth m near trump private student loan rate july loan rate decrease mortgage rate today july rate hold firm this week personal loan rate july rate move down today mortgage refinance rate july rate decrease aca marketplace health plan enrollee could face subsidy cliff in here how to avoid it synchrony report second quarter result company also announces quarterly common stock dividend of per share top breakout altcoins poised for explosive growth sui ofnt and kaspa danaher report second quarter result quavo fraud dispute secures million growth investment from spectrum equity truman warns scare buyer ngex mineral to spinout royalty to shareholder huya inc to report second quarter financial result on tuesday august city move towards major commercial development on city southside exyte receives facility diamond partner award for catl battery cell gigafactory in germany openai exec fidji simo on how ai can turn idea into income apple iphone series major display upgrade teased ahead of september launch synchrony report second quarter result company also announces quarterly common stock dividend of per share people bancorp inc announces second quarter result liquid death and sheetz add a side of mayhem to lunch r e p e a t medium advisory minister dabrusin to make an announcement in vancouver to highlight the importance of empowering young canadian to address climate change and support a healthy environment minimus announces availability of secure container image on iron bank accelerating fedramp and government compliance nano lab expands bnb strategic reserve to u million and announces plan to upgrade bnb reserve strategy best samsung smartphones under r for vloggers osisko metal gaspe copper project intersects metre averaging cu the vistria group promotes max vielmetti shane maloney michael shaid and steve upham to principal amplify energy announces strategic update and ceo transition avoid penalty submit your q tax form before the looming deadline with ezpaycheck unlocking federal land offer limited relief in nation housing shortage realtorcom find titan mining close landmark usm credit agreement with u exim camtek to report second quarter financial result on tuesday august progress software announces amended credit facility kucoin resume brl fiat balance transaction strengthening commitment to brazilian trader unlocking federal land offer limited relief in nation housing shortage realtorcom find massachusetts public college and university and ecampuscom announce online bookstore partnership clear to announce second quarter financial result on august airline face climate reality check with green jet fuel hunt scanlon venture serf a buyside ma advisor to egon zehnder in acquisition of the prince houston group liberty gold strengthens development team to advance black pine through feasibility and permitting titan mining close landmark usm credit agreement with u exim simmons first national corporation announces pricing of offering of common stock exyte receives facility diamond partner award for catl battery cell gigafactory in germany severe obesity in u child ha more than tripled since desjardins investment make change to it mutual fund lineup people bancorp inc announces second quarter result danaher report second quarter result rokmaster provides update for fieldwork on nechako project simmons first national corporation announces pricing of offering of common stock qinetiq u appoints scott mcgleish a chief growth officer how to grow your business by constantly raising the bar eshallgo inc nasdaq ehgo expands into enterpriselevel ai solution targeting china rapidly growing intelligent office market toncoin price prediction toncoin price could reach this year but so could this altcoin what should you buy global defense contractor burke product selects datavault ai for enhanced national defense and aerospace technology contract the vistria group promotes max vielmetti shane maloney michael shaid and steve upham to principal progress software announces amended credit facility clear to announce second quarter financial result on august money market interest rate today july earn up to ai is spreading on main street and transforming the american small business bloom energy nysebe given new price target at susquehanna csw industrials inc nysecsw share sold by mutual of america capital management llc jpmorgan chase co issue positive forecast for ford motor nysef stock price cohen steer nysecns price target lowered to at bank of america mutual of america capital management llc ha stake in huntington ingalls industry inc nysehii decker outdoor nysedeck price target cut to by analyst at well fargo company under armour inc nyseuaa stock position lessened by mutual of america capital management llc mutual of america capital management llc lower position in badger meter inc nysebmi mutual of america capital management llc acquires new share in element solution inc nyseesi mutual of america capital management llc sell share of sps commerce inc nasdaqspsc the arena group nysearen now covered by lake street capital gate industrial corporation plc nysegtes stake reduced by mutual of america capital management llc compass group loncpg receives buy rating from shore capital capricorn energy cne buy rating reaffirmed at shore capital hvivo lonhvo given buy rating at shore capital fiserv nysefi rating increased to neutral at monness crespi hardt magnolia oil gas corp nysemgy share sold by mutual of america capital management llc comerica nysecma given a price target at da davidson dana nysedan price target cut to by analyst at jpmorgan chase co stifel nicolaus issue positive forecast for flowserve nysefls stock price fiserv nysefi price target lowered to at morgan stanley cushman wakefield nysecwk now covered by citizen jmp arbuthnot banking group arbb house stock rating reaffirmed at shore capital mutual of america capital management llc sell share of franklin resource inc nyseben cummins nysecmi price target raised to at barclays trump blast massie a the worst republican congressman and say he seeking a challenger to support canada tariff could add to the cost of building a home by report warns it just wildly misleading why the administration latest allegation about the russia investigation dont add up chris spear america trucker need a place to park before it too late the trump administration claim no one ha died due to u aid cut our trip to afghanistan suggests otherwise editorial a new cta day for edgewater and uptown but what about broadway sheldon h jacobson oldschool manufacturing is gone and is not coming back editorial cook county democrat got it right by staying neutral daniel depetris donald trump and benjamin netanyahu dont have the same goal pod of war like a zombie apocalypse trump budget cut stir fear of frightening pipeline mishap what to do if your bos is playing favoritesand youre not the favorite my view pie define connection in this extended family state cop quietly tag thousand a gang member and feed their name to ice federal government to take over atlantic city housing authority van drew announces some virginian worry about whether theyll lose health coverage tulsi gabbards dangerous disinformation campaign against america the latest from the middle east a israel ramp up gaza offensive and starvation crisis spiral what is the slang term for the final card turned over in texas hold em venus williams make winning return to tennis aged the man who brought you trump bible in school is truly back at it again kellanova nysek receives average recommendation of hold from analyst juror in the james craig trial can ask witness question how it work and what theyve asked raising a family is not easy here innout owner fleeing california moving east realty income corporation nyseo given consensus rating of hold by analyst analyst set bellring brand inc nysebrbr pt at cohu inc nasdaqcohu given average recommendation of moderate buy by analyst in a group of spy had a targetand a plan it turned into one of the most sensationally botched crime of the century facility where greensboro man went missing is unlicensed state regulator say analyst set salesforce inc nysecrm price target at trump commerce secretary dish on president bizarre bedtime call trump admin ditch unesco say it a hotbed of woke and antiisrael propaganda phonics and reading skill take a big step forward other view phonics and reading skill take a big step forward other view epic city replacing old glory with the crescent moon garbarino poised to lead house homeland security panel after green exit going to the beach look very different in britain ct bottle bill now tackle redemption fraud but is it happening only one scandal could really hurt trump and make him more dangerous than ever there no safe place for legal immigrant in trump america burke county get grant for trail including morganton greenway fonta flora trail justice dept anticipates ghislaine maxwell meeting in coming day amid epstein fallout are trump supporter showing sign of moving on from the jeffrey epstein story letter hero too can come from somewhere else letter hero too can come from somewhere else ssc technology ssnc to release quarterly earnings on wednesday heavy metal legend ozzy osbourne dy at on epstein democrat go low from grassroots competition to global industry the evolution of esports brokerage set nvidia corporation nasdaqnvda target price at saia inc nasdaqsaia given average rating of moderate buy by brokerage brokerage set definitive healthcare corp nasdaqdh target price at student protest after military jet crash kill mostly child tennessee titan qb will levi will have seasonending shoulder surgery coldplay fan reveals if shes made money from that astronomer ceo viral concert video battenfeld where is josh kraft path to victory letter thankful for decision blocking stonegate project letter thankful for decision blocking stonegate project how many shark attack in florida what to do to reduce chance of being bitten by shark trump labor department proposes more than rule change in a push to deregulate workplace boston pushing for more everett stadium mitigation heightening tension with kraft who will succeed fed chair jerome powell here are the top contender opposite of reality critic torch karoline leavitts bold new claim about trump faq here what would happen if trump actually try to fire powell trump a foil disappears from andrew cuomos rebooted new york mayoral bid liz peek trump could secure his place in american history with this bold immigration reform strategy opinion celebrating year of aviation innovation and investing in connecticut future fifth ave from boulevard to driveway and back mamdani must back yeshiva standard stripping npr funding shuts down voice of dissent dont withhold promised adult education fund there a growing sentiment gap between rich and poor american scoop the overlooked data trump economist see predicting a boom ahead trump forever grudge tie epstein file to russia hoax it crucial for u to become energy independent spartanburg county council district gop race to succeed britt defined by county rapid growth mini tds is psychosis defined a an irrational fear search for missing yearold greensboro man continues family offer reward criticizes police option to pay for college just got worse stryker corporation nysesyk receives average recommendation of moderate buy from analyst how schumer is preparing for the fall funding showdown commentary trump problem determining putin desire for peace commentary supreme court silent opinion undermine it legitimacy recording dream is exhausting why male corporate leader and billionaire need financial therapy ai grows up without rule san francisco to ban homeless people from living in rv with new parking limit editorial stephen colbert roasted boss and took risk bravo editorial bondis mishandling of epstein case is her mo letter will biomass be our dirty little secret letter will biomass be our dirty little secret bryan kohberger ha another chance to explain why he killed idaho college student but expert say he unlikely to speak judge partially block trump admin push to slash planned parenthood funding concern grow over australia toxic algae bloom commentary why neophyte mayor need steady hand on the wheel commentary if haiti ha become more violent why end haitian temporary protected status in the u editorial make it clear that floridian hate the idea of a blackbear massacre commentary kid in camp nope got a summer schedule nope cue the mom guilt maine transgender student deserve better letter out to get golden republican ensnare collins steve collins broadway box office sunset blvd end run on million high umgs downtown deal european commission launch indepth second phase of investigation accenture expands in malaysia with aristal purchase ap business summarybrief at pm edt the government wa once a steady partner for nonprofit thats changing the government wa once a steady partner for nonprofit thats changing one transporter two critial job align production system delivers dualuse sav verizons plan to make switching harder for customer hit a snag why this analog semiconductor stock is a top undervalued pick abbvie inc nyseabbv share acquired by obermeyer wealth partner vsm wealth advisory llc increase holding in vanguard value etf nysearcavtv american exceptionalism isnt dead foreign investor flow into u show obermeyer wealth partner trim stake in lowes company inc nyselow full circle benchmark retail diesel price is about where it wa a year ago flashback bill maher admitted stephen colbert feed machine bank julius baer co ltd zurich trim stock position in accenture plc nyseacn bank julius baer co ltd zurich sell share of exxon mobil corporation nysexom united bank sell share of eli lilly and company nyselly procter gamble company the nysepg holding raised by united bank alecta tjanstepension omsesidigt ha million holding in bank of america corporation nysebac exxon mobil corporation nysexom share acquired by flynn zito capital management llc nelson capital management llc increase holding in verizon communication inc nysevz alphabet inc nasdaqgoog share sold by flynn zito capital management llc vsm wealth advisory llc acquires new position in alphabet inc nasdaqgoog highlander partner lp invests million in invesco sp equal weight etf nysearcarsp bank julius baer co ltd zurich raise holding in procter gamble company the nysepg highlander partner lp make new million investment in invesco sp equal weight etf nysearcarsp caterpillar inc nysecat share sold by bank julius baer co ltd zurich flynn zito capital management llc trim position in alphabet inc nasdaqgoog amazon is selling a spacious drawer dresser for only and it the perfect storage solution how to use deepseek for stock market analysis keycorp raise nii and loan outlook after solid quarter ignite planner llc acquires share of walmart inc nysewmt hartford fund management co llc acquires share of american express company nyseaxp ignite planner llc ha million stake in walmart inc nysewmt ifp advisor inc increase stake in amgen inc nasdaqamgn ifp advisor inc raise stock position in general dynamic corporation nysegd hartford fund management co llc sell share of general dynamic corporation nysegd rideup trailrush a fullsuspension middrive emtb built for the trail columbia issue expulsion pull degree over propalestinian protest a pressure from trump administration mount dr horton pulte group share move higher after strong earnings independent bank corporation announces quarterly cash dividend on common stock independent bank corporation announces quarterly cash dividend on common stock att inc nyset share purchased by proficio capital partner llc financial sense advisor inc reduces stock position in att inc nyset yield ease ahead of federal reserve meeting should you buy the postearnings dip in lockheed martin stock former sandoz inc ceo don degolyer lead pharmatech startup ategenos a company prepares public launch lori calvasina theory that tariff wouldnt have an impact are not true guidance capital inc ha million stock holding in freeportmcmoran inc nysefcx klp kapitalforvaltning a lower holding in m company nysemmm klp kapitalforvaltning a ha million position in the sherwinwilliams company nyseshw klp kapitalforvaltning a grows stock holding in analog device inc nasdaqadi lockheed record b in loss mostly linked to continued strife on classified aero program breaking defense u the philippine struck trade military deal trump say openai seek additional capital from investor a part of it billion round wired pier investment partner llc make new investment in salesforce inc nysecrm klp kapitalforvaltning a raise holding in cme group inc nasdaqcme bristol myers squibb company nysebmy share purchased by klp kapitalforvaltning a klp kapitalforvaltning a acquires share of salesforce inc nysecrm salesforce inc nysecrm share purchased by klp kapitalforvaltning a klp kapitalforvaltning a buy share of cme group inc nasdaqcme klp kapitalforvaltning a ha million stock holding in bristol myers squibb company nysebmy klp kapitalforvaltning a purchase share of automatic data processing inc nasdaqadp gm say ev are it north star a legacy automaker chase tesla cnbc indonesia to cut tariff ontariff barrier in u trade deal trump official say amy g injury firm welcome multilingual legal expert to support serbianspeaking client jpmorgan chase co issue positive forecast for quantum service nysepwr stock price universal beteiligungs und servicegesellschaft mbh buy share of intuit in
'@

$i22Text = @'
new commerce split financial result to may prospect of higher indian sugar export weigh on price meta launch new teen safety feature remove account that sexualize child meta launch new teen safety feature remove account that sexualize child new commerce split financial result to may skydance pledge to trump fcc itll eliminate dei install ombudsman to root out bias at cbs news how trump ai plan would cut regulation and speed up adoption stock jump on trump trade deal deluge ahead of big tech earnings whats moving market wednesday cocoa price jump on fear of tighter supply qfscoin launch the most profitable platform for yielding crypto in designed for global investor to earn btc ltc and doge coin united airline cancel all flight to sunny destination offer refund coffee price underpinned by frost concern in brazil procter gamble company the nysepg share bought by j w coon advisor llc poplar forest capital llc sell share of international business machine corporation nyseibm international business machine corporation nyseibm position increased by j w coon advisor llc j w coon advisor llc cut stake in oracle corporation nyseorcl burney co reduces stock position in procter gamble company the nysepg poplar forest capital llc lower position in international business machine corporation nyseibm j w coon advisor llc raise position in international business machine corporation nyseibm hasbro relies on magic to offset impact of tariff prime dividend corp financial result to may best biotech stock to buy according to billionaire steve cohen prime dividend corp financial result to may u automaker say trump tariff deal with japan put them at a disadvantage u automaker say trump tariff deal with japan put them at a disadvantage u automaker say trump tariff deal with japan put them at a disadvantage u automaker say trump tariff deal with japan put them at a disadvantage u automaker say trump tariff deal with japan put them at a disadvantage u automaker say trump tariff deal with japan put them at a disadvantage u automaker say trump tariff deal with japan put them at a disadvantage u automaker say trump tariff deal with japan put them at a disadvantage u automaker say trump tariff deal with japan put them at a disadvantage u automaker say trump tariff deal with japan put them at a disadvantage u automaker say trump tariff deal with japan put them at a disadvantage u automaker say trump tariff deal with japan put them at a disadvantage u automaker say trump tariff deal with japan put them at a disadvantage u automaker say trump tariff deal with japan put them at a disadvantage u automaker say trump tariff deal with japan put them at a disadvantage u automaker say trump tariff deal with japan put them at a disadvantage summit wealth partner llc ha holding in the home depot inc nysehd summit wealth partner llc grows holding in tesla inc nasdaqtsla wedbush security inc purchase share of the home depot inc nysehd wedbush security inc ha million stock position in the home depot inc nysehd egerton capital uk llp raise stock position in eli lilly and company nyselly republican skittish over epstein vote close u house early trump ai action plan call for deregulation chatbot free speech merz macron huddle over nato and euus trade row firsttime homebuyer guide every step detailed from budgeting to close firsttime homebuyer guide every step detailed from budgeting to close firsttime homebuyer guide every step detailed from budgeting to close firsttime homebuyer guide every step detailed from budgeting to close firsttime homebuyer guide every step detailed from budgeting to close first student awarded nstas go yellow go green award for leadership in advancing zeroemission school transportation best ai phone under r in historic order for usbuilt lng carrier could test new rule iran hit by power protest a oil export surge to million bpd first student awarded nstas go yellow go green award for leadership in advancing zeroemission school transportation forget old meme coin this newcomer below will skyrocket and outperform dogecoin doge walmart is selling a soft summer dress for that flow well is the trump account enough how to maximize your child financial future canadian banc corp financial result to may desktop search traffic see surge in ai chatbot usage north american financial split corp financial result to may fhlbank chicago announces q financial highlight food price follow inflation uptick beef hit record high north american financial split corp financial result to may canadian banc corp financial result to may grant cardone slam california plan to buy fireravaged palisade for lowincome housing say it not just a conspiracy year bond auction see robust demand best logistics stock to buy according to hedge fund delta is using ai for flight pricingheres how to save with credit card landmark wealth management inc acquires share of vanguard sp etf nysearcavoo vanguard sp etf nysearcavoo share acquired by wakefield asset management lllp ames national corp boost stake in vanguard sp etf nysearcavoo canadian life company split corp financial result to may canadian life company split corp financial result to may how leader undermine their own authority macron sue candace owen for defamation for claiming his wife is a man u financial split corp financial result to may jim cramer describes positive and negative force that are shaping market action judge denies trump admin bid to unseal jeffrey epstein grand jury transcript in florida david letterman had something to say after the stephen colbert cancellation iqvia holding iqv snap day loss on higherthanexpected earnings dr horton dhi jump after earnings beat new doubleday diner opening in ballston spa in old route diner spot webull bull soar on a stronger crypto adoption bolster business were raising our price target on ge vernova by after blowout earnings report cnbc were raising our price target on ge vernova by after blowout earnings report pultegroup phm snap day loss on buy reco higher price target hershey and other chocolate maker hike price a cocoa remains near record high hershey and other chocolate maker hike price a cocoa remains near record high hershey and other chocolate maker hike price a cocoa remains near record high hershey and other chocolate maker hike price a cocoa remains near record high cipher mining cifr jump alongside bitcoin ahead of q result trend of equity market is still higher but investor are fighting it piper sandlers craig johnson the gsa just proposed big change that could affect small business contractor dividend split corp financial result to may lucid group lcid soar a vehicle gain access to tesla charger nio nio jump on new vehicle launch quantumscape q gain ground after loss up ahead of q earnings scalise house will deal with epstein issue after august recess if needed skydance pledge to trump fcc itll eliminate dei install ombudsman to root out bias at cbs news candace owen sued by french president over trans claim about wife bryan kohberger speaks just three word when given chance to explain why he killed four idaho college student roy cooper to jump into north carolina senate race monday j w coon advisor llc sell share of johnson johnson nysejnj boyd watterson asset management llc oh buy share of johnson johnson nysejnj the mindful minute by jon heydenreich pick me ups u automaker say trump tariff deal with japan put them at a disadvantage u automaker say trump tariff deal with japan put them at a disadvantage u automaker say trump tariff deal with japan put them at a disadvantage u automaker say trump tariff deal with japan put them at a disadvantage u automaker say trump tariff deal with japan put them at a disadvantage u automaker say trump tariff deal with japan put them at a disadvantage u automaker say trump tariff deal with japan put them at a disadvantage u automaker say trump tariff deal with japan put them at a disadvantage u automaker say trump tariff deal with japan put them at a disadvantage u automaker say trump tariff deal with japan put them at a disadvantage u automaker say trump tariff deal with japan put them at a disadvantage egerton capital uk llp raise holding in eli lilly and company nyselly release the damn file outgoing gop senator slam trump over epstein how can we give feedback that people actually use the ethical life podcast how can we give feedback that people actually use the ethical life podcast how can we give feedback that people actually use the ethical life podcast how can we give feedback that people actually use the ethical life podcast how can we give feedback that people actually use the ethical life podcast how can we give feedback that people actually use the ethical life podcast how can we give feedback that people actually use the ethical life podcast how can we give feedback that people actually use the ethical life podcast how can we give feedback that people actually use the ethical life podcast how can we give feedback that people actually use the ethical life podcast trump administration violated the law by withholding some head start fund congressional watchdog find live now senate committee examines u diplomatic strategy in the middle east thailand close border crossing with cambodia and recall ambassador a tension flare notis mcconarty edward acquires share of amgen inc nasdaqamgn wakefield asset management lllp acquires share of vanguard sp etf nysearcavoo ames national corp buy share of visa inc nysev spdr gold share nysearcagld share sold by landmark wealth management inc speaker johnson vow appropriate action in house if epstein file arent released speaker johnson vow appropriate action in house if epstein file arent released speaker johnson vow appropriate action in house if epstein file arent released david letterman had something to say after the stephen colbert cancellation speaker johnson vow appropriate action in house if epstein file arent released update trump wa told he is in epstein file wall street journal report update trump wa told he is in epstein file wall street journal report speaker johnson vow appropriate action in house if epstein file arent released speaker johnson vow appropriate action in house if epstein file arent released live now senate committee hold hearing on innovation and china industrial espionage u agency warn of surging interlock ransomware attack targeting healthcare business white house blast wall street journal report that trump is in epstein file trump wa told he is in epstein file wall street journal report wespac advisor socal llc grows position in crowdstrike nasdaqcrwd wespac advisor socal llc invests in costco wholesale corporation nasdaqcost replay house armed service committee hold hearing on reforming defense acquisition forbes j m co llp buy share of pepsico inc nasdaqpep gabbard leaving obama admin criminal implication to doj fbi in russia hoax contrived narrative pa gov josh shapiro call supporter of ban on transgender athlete extremist legislator decline to say if hed veto trygve olson cartoon state fair ride susan collins will oppose making trump former lawyer a federal judge jpl wealth management llc ha million stock holding in caterpillar inc nysecat caterpillar inc nysecat share sold by close asset management ltd jpl wealth management llc take position in ishares core sp smallcap etf nysearcaijr youre not imagining it there ha been a bumper crop of special election in minnesota legislature youre not imagining it there ha been a bumper crop of special election in minnesota legislature what trump and pam bondi are doing in new jersey is a bigger deal than you think live now senate judiciary committee hold hearing on ending discrimination and preference steel plant in vances hometown trade clean future for more coal share in palantir technology inc nasdaqpltr purchased by garrison bradford associate inc live now house oversight committee hold hearing on mail theft and crime no there is not a tea bag killer in detroit police say hareholder alert the ma class action firm announces an investigation of susquehanna community financial inc otcmkts sqcf merck co inc nysemrk stock holding decreased by poplar forest capital llc rep robert bresnahan jr purchase share of abbvie inc nyseabbv boyd watterson asset management llc oh acquires share of abbvie inc nyseabbv trump ally lindell win appeal in lawsuit over million election contest mike lawler pass on ny governor run seek reelection in swing seat summit wealth partner llc lower stock position in johnson johnson nysejnj should college athlete be employee house panel say no grand fork to host it firstever pride parade in september bryan kohberger sentenced to life for idaho student killing open retribution alleged a san diegos doj office hit by resignation demotion letter putting tariff on brazil a bad move what we know about the epstein file nebraska sen pete ricketts set to introduce bill streamlining veteran health care famed south african activist get back into fight for hiv med after usaid cut live now house oversight and government reform hearing addressing the crisis of unaccompanied alien childrenpart i doctor pleads guilty to distributing ketamine tied to matthew perry death pelham nh man charged after cat found abandoned in crate outside pet center one year reopening timeline for veteran home in enterprise scoop cooper plan to enter north carolina senate race next week roundview capital llc increase stock position in netflix inc nasdaqnflx salesforce inc nysecrm share sold by ima advisory service inc idaho student murderer bryan kohberger in court to be sentenced alphabet inc nasdaqgoog share sold by close asset management ltd city national bank of florida msd purchase share of alphabet inc nasdaqgoog crowdstrike nasdaqcrwd share purchased by legato capital management llc supreme court permit firing at consumer product safety agency judge sentence bryan kohberger to life in prison in idaho quadruplemurder case rebuilding syria opportunity and challenge of postwar reconstruction new round of russiaukraine talk start in turkey day after trump gave moscow deadline to make peace
'@

# Row 21
$ws.Range("A21").Value = 45861
$ws.Range("A21").NumberFormat = $ws.Range("A20").NumberFormat
$ws.Range("B21").Value = 6358.91015625
$ws.Range("C21").Value = 6360.64013671875
$ws.Range("D21").Value = 6317.490234375
$ws.Range("E21").Value = 6326.89990234375
$ws.Range("F21").Value = 5642510000
$ws.Range("G21").Value = 0.0078118869515191
$ws.Range("H21").Value = 1
$ws.Range("I21").Value = $i21Text

# Row 22
$ws.Range("A22").Value = 45862
$ws.Range("A22").NumberFormat = $ws.Range("A20").NumberFormat
$ws.Range("B22").Value = 6363.35009765625
$ws.Range("C22").Value = 6381.31005859375
$ws.Range("D22").Value = 6360.56982421875
$ws.Range("E22").Value = 6368.60009765625
$ws.Range("F22").Value = 3364122000
$ws.Range("G22").Value = 0.0006982236416543
$ws.Range("H22").Value = 1
$ws.Range("I22").Value = $i22Text
